# Refactor Message classes into a Message package
#
# Target paragraph (the one ending "...will describe Java classes.") gets:
#   1. New text prepended in front of it (ending in a lower-case "a" that
#      continues into the original "ll communication..." text).
#   2. New text appended after it (" All data will be wrapped in...").
#   3. A single trailing space appended after the existing _GoBack bookmark.
#
# The original run "All communication is ... Java classes." needs to be
# split after its first character ("A" / "ll communication...") so the
# first run ends in "...for agents, a" and a second, distinct run begins
# with "ll communication...".

$d = $word.ActiveDocument

# Locate the target paragraph (the last real paragraph of the document,
# which contains the JADE INFORM sentence) by searching for its text.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.StartsWith("All communication is implemented")) {
        $target = $cand
        break
    }
}
if ($target -eq $null) {
    throw "Could not locate target paragraph"
}

$quoteLeft = [char]0x201C
$quoteRight = [char]0x201D
$apostrophe = [char]0x2019

$run1Text = "Another goal of the communication protocol is to avoid re-inventing the wheel. JADE provides a transport layer to send Java objects as messages between JADE agents. Thus, we use Jade" + $apostrophe + "s Agent Identifier (AID) identification system for agents, a"
$run3Text = " All data will be wrapped in " + $quoteLeft + "Message" + $quoteRight + " objects in order to separate message processing logic from data model. The data model should function independently of the messaging layer used to move it."
$run4Text = " "

# --- Step 1: delete the leading "A" of "All communication..." -------------
# This leaves the rest of the original sentence ("ll communication ... Java
# classes.") as its own, untouched run, and leaves the existing
# bookmarkStart/bookmarkEnd (_GoBack) exactly where they were.
$full = $target.Range
$rA = $d.Range($full.Start, $full.Start + 1)
$rA.Delete()

# --- Step 2: prepend the new "Another goal ..." run ------------------------
# A zero-length Range positioned at the very start of the paragraph, passed
# to InsertXML, inserts its content as a brand-new run *before* the
# existing run rather than merging into it.
$rStart = $d.Range($full.Start, $full.Start)
$xml1 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>' + $run1Text + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rStart.InsertXML($xml1)

# --- Step 3: append the new trailing text after the whole paragraph -------
# InsertAfter on a Range collapsed at the paragraph's end appends plain
# text after the existing bookmarkEnd, merging formatting-compatible runs -
# that's fine here since we immediately re-split it in step 4.
# Re-fetch the paragraph object/range after the prior edits.
$paraRange = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.StartsWith("Another goal of the communication protocol")) {
        $paraRange = $cand.Range
        break
    }
}
$endR = $d.Range($paraRange.End, $paraRange.End)
$marker = "|||SPLITMARKER|||"
$endR.InsertAfter($run3Text + $marker + $run4Text)

# --- Step 4: move the _GoBack bookmark to sit between run3 and run4 -------
# Re-read the paragraph text/position, find the marker, and re-add the
# bookmark there; Bookmarks.Add repositions an existing bookmark of the
# same name. This also has the side effect of splitting the run that
# currently holds "<run3Text><marker><run4Text>" into two runs at that
# point, which is exactly the run boundary the target document has.
$paraRange = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.StartsWith("Another goal of the communication protocol")) {
        $paraRange = $cand.Range
        break
    }
}
$paraText = $paraRange.Text
$splitIdx = $paraText.IndexOf($marker)
$bmPos = $paraRange.Start + $splitIdx
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# --- Step 5: strip the temporary marker text out ---------------------------
$paraRange = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.StartsWith("Another goal of the communication protocol")) {
        $paraRange = $cand.Range
        break
    }
}
$paraText = $paraRange.Text
$markerIdx = $paraText.IndexOf($marker)
$markerRange = $d.Range($paraRange.Start + $markerIdx, $paraRange.Start + $markerIdx + $marker.Length)
$markerRange.Delete()
